$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (Adaptive6 / job 768) and row 8 (Allium / job 787) from the original sheet.
# Deleting row 5 first causes row 8 to shift up to row 7, so delete that one next.
$ws.Rows("5").Delete()
$ws.Rows("7").Delete()
